$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 81.25
$ws.Range("I8").Value = 81.25
$ws.Range("K8").Value = 243.75
$ws.Range("M8").Value = -104.75
$ws.Range("H21").Value = 22379.25
$ws.Range("I21").Value = 9258.5
$ws.Range("K21").Value = 9258.5
$ws.Range("M21").Value = -8790.5
$ws.Range("H23").Value = 22379.25
$ws.Range("I23").Value = 9258.5
$ws.Range("K23").Value = 9258.5
$ws.Range("M23").Value = -9024.5
$ws.Range("H38").Value = 640
$ws.Range("I38").Value = 640
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1920
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1548
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 3780
$ws.Range("I58").Value = 3780
$ws.Range("K58").Value = 11340
$ws.Range("M58").Value = -11190
$ws.Range("H69").Value = 46106.918
$ws.Range("I69").Value = 2891.1428
$ws.Range("J69").Value = 106609
$ws.Range("K69").Value = 8673.4284
$ws.Range("L69").Value = 319827
$ws.Range("M69").Value = -7799.428400000001
$ws.Range("N69").Value = -321575
$ws.Range("H72").Value = 46106.918
$ws.Range("I72").Value = 2891.1428
$ws.Range("J72").Value = 106609
$ws.Range("K72").Value = 26020.2852
$ws.Range("L72").Value = 959481
$ws.Range("M72").Value = -21652.2852
$ws.Range("N72").Value = -968217
$ws.Range("H74").Value = 11852.5
$ws.Range("I74").Value = 12058.444
$ws.Range("J74").Value = 9999
$ws.Range("K74").Value = 12058.444
$ws.Range("L74").Value = 9999
$ws.Range("M74").Value = -11122.444
$ws.Range("N74").Value = -11871
$ws.Range("H77").Value = 11852.5
$ws.Range("I77").Value = 12058.444
$ws.Range("J77").Value = 9999
$ws.Range("K77").Value = 60292.22
$ws.Range("L77").Value = 49995
$ws.Range("M77").Value = -55612.22
$ws.Range("N77").Value = -59355
$ws.Range("H87").Value = 59552.332
$ws.Range("J87").Value = 59552.332
$ws.Range("L87").Value = 59552.332
$ws.Range("N87").Value = -62048.332
$ws.Range("H90").Value = 59552.332
$ws.Range("J90").Value = 59552.332
$ws.Range("L90").Value = 178656.996
$ws.Range("N90").Value = -191136.996
$ws.Range("H135").Value = 380.14285
$ws.Range("I135").Value = 371.4
$ws.Range("J135").Value = 555
$ws.Range("K135").Value = 3342.6
$ws.Range("L135").Value = 4995
$ws.Range("M135").Value = -807.5999999999999
$ws.Range("N135").Value = -10065
$ws.Range("H137").Value = 1754.1111
$ws.Range("J137").Value = 2121.4
$ws.Range("L137").Value = 6364.200000000001
$ws.Range("N137").Value = -11464.2
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8598.031
$ws.Range("I32").Value = 8240.931
$ws.Range("J32").Value = 12050
$ws.Range("K32").Value = 8240.931
$ws.Range("L32").Value = 12050
$ws.Range("M32").Value = -7953.931
$ws.Range("N32").Value = -12624
$ws.Range("H45").Value = 1886
$ws.Range("I45").Value = 1810.7142
$ws.Range("J45").Value = 2149.5
$ws.Range("K45").Value = 1810.7142
$ws.Range("L45").Value = 2149.5
$ws.Range("M45").Value = -1433.7142
$ws.Range("N45").Value = -2903.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 20500
$ws.Range("I35").Value = 5750
$ws.Range("K35").Value = 5750
$ws.Range("M35").Value = -5456
$ws.Range("H125").Value = 15000
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 15000
$ws.Range("N125").Value = -19920
$ws.Range("H134").Value = 2672.4
$ws.Range("I134").Value = 2225
$ws.Range("K134").Value = 6675
$ws.Range("M134").Value = -4140
$ws.Range("H141").Value = 272999.34
$ws.Range("J141").Value = 272999.34
$ws.Range("L141").Value = 272999.34
$ws.Range("N141").Value = -283359.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2636.625
$ws.Range("I137").Value = 2274.75
$ws.Range("K137").Value = 6824.25
$ws.Range("M137").Value = -1724.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18668
$ws.Range("I5").Value = 18668
$ws.Range("K5").Value = 18668
$ws.Range("M5").Value = -18556
$ws.Range("H70").Value = 2999
$ws.Range("I70").Value = 2999
$ws.Range("K70").Value = 2999
$ws.Range("M70").Value = -2729
$ws.Range("H73").Value = 2999
$ws.Range("I73").Value = 2999
$ws.Range("K73").Value = 2999
$ws.Range("M73").Value = -2063
$ws.Range("H107").Value = 2055.3333
$ws.Range("I107").Value = 3152.6
$ws.Range("K107").Value = 3152.6
$ws.Range("M107").Value = -1232.6
$ws.Range("H126").Value = 2465.6667
$ws.Range("I126").Value = 2465.6667
$ws.Range("K126").Value = 7397.000100000001
$ws.Range("M126").Value = -4927.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1777.7391
$ws.Range("I93").Value = 1531.2858
$ws.Range("J93").Value = 2161.111
$ws.Range("K93").Value = 1531.2858
$ws.Range("L93").Value = 2161.111
$ws.Range("M93").Value = -283.2858000000001
$ws.Range("N93").Value = -4657.111
$ws.Range("H136").Value = 3760.4443
$ws.Range("I136").Value = 3760.4443
$ws.Range("K136").Value = 11281.3329
$ws.Range("M136").Value = -8731.332900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 632.5
$ws.Range("J107").Value = 671.5714
$ws.Range("L107").Value = 2014.7142
$ws.Range("N107").Value = -5854.7142
$ws.Range("H132").Value = 4998.6665
$ws.Range("J132").Value = 4998.6665
$ws.Range("L132").Value = 14995.9995
$ws.Range("N132").Value = -20055.9995
$ws.Range("H136").Value = 5574.15
$ws.Range("I136").Value = 4749.2144
$ws.Range("J136").Value = 7499
$ws.Range("K136").Value = 14247.6432
$ws.Range("L136").Value = 22497
$ws.Range("M136").Value = -11697.6432
$ws.Range("N136").Value = -27597
